# Update build timestamp in version strings from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on " + $newStamp + ")"

$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Uvalnaya Coal Mine, Russia, M1520, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$newVersionText = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on " + $newStamp + ")"

for ($row = 2; $row -le 17; $row++) {
    $wsData.Range("S" + $row).Value = $newVersionText
}
